# Append the new March-2023 data row (row 76) to Sheet1, matching the
# source INDEC dataset's monthly refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 76

$ws.Cells.Item($row, 1).Value  = 7.7                 # A76 Nivel general
$ws.Cells.Item($row, 2).Value  = 9.300000000000001   # B76 Alimentos y bebidas
$ws.Cells.Item($row, 3).Value  = 8.300000000000001   # C76 Bebidas alcohólicas y tabaco
$ws.Cells.Item($row, 4).Value  = 9.4                 # D76 Prendias y Calzado
$ws.Cells.Item($row, 5).Value  = 6.5                 # E76 Vivienda Agua y Elec
$ws.Cells.Item($row, 6).Value  = 5.8                 # F76 Equip y Mant del Hogar
$ws.Cells.Item($row, 7).Value  = 5.7                 # G76 Salud
$ws.Cells.Item($row, 8).Value  = 5.3                 # H76 Transporte
$ws.Cells.Item($row, 9).Value  = 1.9                 # I76 Comunicación
$ws.Cells.Item($row, 10).Value = 4.4                 # J76 Recreación y cultura
$ws.Cells.Item($row, 11).Value = 29.1                # K76 Educación
$ws.Cells.Item($row, 12).Value = 7.9                 # L76 Restaurantes y hoteles
$ws.Cells.Item($row, 13).Value = 6.3                 # M76 Bienes y servicios varios

# N76 periodos - date serial 44986 == 2023-03-01. Set the value, then copy
# the date number-format down from the cell above so it reuses the same
# yyyy-mm-dd style already in the workbook instead of minting a new one.
$ws.Cells.Item($row, 14).Value = 44986
$ws.Range("N75").Copy()
$ws.Cells.Item($row, 14).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item($row, 15).Value = 2023                # O76 year
$ws.Cells.Item($row, 16).Value = 3                   # P76 Mes
$ws.Cells.Item($row, 17).Value = "Mar"                # Q76 month

$wb.Save()
